# Add a new food-log entry as row 12, mirroring the existing rows (1,2,11)
# which store every value (including numbers/dates) as plain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 12
$values = @("07/07/2025", "3 eggs omlette", "234", "18", "2", "17", "561")

# Force the range to Text format first so Excel doesn't reinterpret the
# date-looking / numeric-looking strings as a date serial or a number.
$rowRange = $ws.Range("A$newRow`:G$newRow")
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = $values[0]
$ws.Cells.Item($newRow, 2).Value = $values[1]
$ws.Cells.Item($newRow, 3).Value = $values[2]
$ws.Cells.Item($newRow, 4).Value = $values[3]
$ws.Cells.Item($newRow, 5).Value = $values[4]
$ws.Cells.Item($newRow, 6).Value = $values[5]
$ws.Cells.Item($newRow, 7).Value = $values[6]

# Reset the style back to the sheet default so the new cells don't carry an
# explicit "Text" number-format style the way the original rows don't either.
$rowRange.Style = "Normal"
